# Commit: "commit on 16th Final"
#
# This script reproduces the data / view-state edits from the target
# commit:
#   - Test_Scenarios!D4 changes from "Yes" to "No"
#   - Keyword_Scenario!C3 and !C4 change from "Application_Submit" to
#     "Application_Submit1"
#   - The active sheet/tab moves from "Test_Scenarios" to
#     "Keyword_Scenario" (tabSelected moves accordingly), and the
#     selected cell on Keyword_Scenario becomes C4 (instead of D7).

$wb = $excel.ActiveWorkbook

# --- Test_Scenarios: D4 "Yes" -> "No" -----------------------------------
$wsTestScenarios = $wb.Worksheets.Item("Test_Scenarios")
$wsTestScenarios.Range("D4").Value = "No"

# --- Keyword_Scenario: C3 & C4 "Application_Submit" -> "Application_Submit1"
$wsKeyword = $wb.Worksheets.Item("Keyword_Scenario")
$wsKeyword.Range("C3").Value = "Application_Submit1"
$wsKeyword.Range("C4").Value = "Application_Submit1"

# --- Make Keyword_Scenario the active/selected sheet, with C4 selected --
$wsKeyword.Activate()
$wsKeyword.Range("C4").Select()
